$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values for the rows that changed
$ws.Range("F4").Value = -13
$ws.Range("F5").Value = -7
$ws.Range("F7").Value = -7
$ws.Range("F8").Value = 0
$ws.Range("F14").Value = -8
$ws.Range("F15").Value = -8
$ws.Range("F16").Value = -5
$ws.Range("F17").Value = -6
$ws.Range("F19").Value = -5
$ws.Range("F20").Value = 3
$ws.Range("F22").Value = -2
$ws.Range("F23").Value = -2
$ws.Range("F26").Value = 4
$ws.Range("F27").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 2
$ws.Range("F37").Value = 0
$ws.Range("F42").Value = -1
